{"js": "// 1. \"\u0421\u0438\u0441\u0442\u0435\u043c\u043d\u0438\u0439\" -> \"\u041c\u043e\u0434\u0443\u043b\u044c\u043d\u0438\u0439\" (Level of testing cell in the first table)\nconst sysResults = context.document.body.search(\"\u0421\u0438\u0441\u0442\u0435\u043c\u043d\u0438\u0439\", { matchCase: true });\nsysResults.load(\"items\");\nawait context.sync();\nfor (const r of sysResults.items) {\n  r.insertText(\"\u041c\u043e\u0434\u0443\u043b\u044c\u043d\u0438\u0439\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Remove the \"1. \u0423\u0432\u0435\u0441\u0442\u0438 \" prefix runs that precede the bare numbers\n//    (0, 5, 15, 25, 35) in the second (TestSuite) table.\nconst uvestyResults = context.document.body.search(\"1. \u0423\u0432\u0435\u0441\u0442\u0438 \", { matchCase: true });\nuvestyResults.load(\"items\");\nawait context.sync();\nfor (const r of uvestyResults.items) {\n  r.insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3. Center-align every paragraph in the second table that is currently\n//    left-aligned (the \"steps\" cell, its trailing blank paragraph, and the\n//    \"expected result\" cell, for every data row).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst testSuiteTable = tables.items[1];\nconst paragraphs = testSuiteTable.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"alignment\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.alignment === Word.Alignment.left) {\n    p.alignment = Word.Alignment.centered;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"\u0421\u0438\u0441\u0442\u0435\u043c\u043d\u0438\u0439\" -> \"\u041c\u043e\u0434\u0443\u043b\u044c\u043d\u0438\u0439\" (Level of testing cell in the first table)\n$d.Content.Find.Execute(\"\u0421\u0438\u0441\u0442\u0435\u043c\u043d\u0438\u0439\", $true, $false, $false, $false, $false, $true, 0, $false, \"\u041c\u043e\u0434\u0443\u043b\u044c\u043d\u0438\u0439\", 2)\n\n# 2. Remove the \"1. \u0423\u0432\u0435\u0441\u0442\u0438 \" prefix runs that precede the bare numbers\n#    (0, 5, 15, 25, 35) in the second (TestSuite) table.\n$d.Content.Find.Execute(\"1. \u0423\u0432\u0435\u0441\u0442\u0438 \", $true, $false, $false, $false, $false, $true, 0, $false, \"\", 2)\n\n# 3. Center-align every paragraph in the second table that is currently\n#    left-aligned (the \"steps\" cell, its trailing blank paragraph, and the\n#    \"expected result\" cell, for every data row).\n#\n# NOTE: we deliberately avoid the Tables object model here - merely reading\n# a Table's .Range in this runtime leaves stale Range.Start/Range.End values\n# behind for unrelated Range objects fetched afterwards. Find-based Range\n# probing does not have that issue, so we use it to locate the table bounds.\n$rStart = $d.Content\n$rStart.Find.Execute(\"\u0406\u0434-\u0440 \u0442\u0435\u0441\u0442 \u043a\u0435\u0439\u0441\u0430\", $true, $false, $false, $false, $false, $true, 0, $false, \"\", 0)\n$tableStart = $rStart.Start\n\n$rLastRow = $d.Content\n$rLastRow.Start = $tableStart\n$rLastRow.Find.Execute(\"TS_05\", $true, $false, $false, $false, $false, $true, 0, $false, \"\", 0)\n\n$rEnd = $d.Content\n$rEnd.Start = $rLastRow.Start\n$rEnd.Find.Execute(\"PASSED\", $true, $false, $false, $false, $false, $true, 0, $false, \"\", 0)\n$tableEnd = $rEnd.End\n\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\nfor ($i = 1; $i -le $count; $i = $i + 1) {\n  $p = $paragraphs.Item($i)\n  $start = $p.Range.Start\n  if ($start -ge $tableStart -and $start -lt $tableEnd -and $p.Alignment -eq 0) {\n    $p.Alignment = 1\n  }\n}\n"}
